$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving plain numeric-looking text must be forced to Text format
# so Excel does not auto-convert them into numbers, matching the original
# inline-string representation.
$textCells = @('D5', 'D6', 'D10', 'D11', 'D15', 'D16', 'D19', 'D21', 'D24', 'D25', 'D26', 'D35', 'D36', 'D40', 'D42', 'D43', 'D44', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.723.79'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.600.01'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '211.36'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '0.514'
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('D10').Value = '19.53'
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('D11').Value = '0.0841'
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '1.823.30'
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').Value = '1.614.75'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('D16').Value = '65.42'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('D17').Value = '26.698.76'
$ws.Range('E17').Value = '  +0.39%  '
$ws.Range('D18').Value = '0.0₃0763'
$ws.Range('E18').Value = '  +4.74%  '
$ws.Range('D19').Value = '210.00'
$ws.Range('E19').Value = '  +1.25%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').Value = '7.18'
$ws.Range('E21').Value = '  +4.87%  '
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').Value = '8.93'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('D25').Value = '143.20'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('E30').Value = '  +3.70%  '
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('D34').Value = '1.290.06'
$ws.Range('E34').Value = '  +0.50%  '
$ws.Range('D35').Value = '0.620'
$ws.Range('E35').Value = '  -4.85%  '
$ws.Range('D36').Value = '2.47'
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('E39').Value = '  +17.14%  '
$ws.Range('D40').Value = '0.829'
$ws.Range('E40').Value = '  -1.11%  '
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.785'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.19'
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('D44').Value = '63.18'
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').Value = '1.736.42'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('D46').Value = '91.42'
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('D47').Value = '1.58'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0104'
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.100'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0508'
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.07%  '
